# Apply "Added runs, scripts and pictures" changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Block 1: rows 45-48  (SOA_formation_24062013_faster)
# ---------------------------------------------------------------------
$ws.Range("A45").Value = "Script"

$ws.Range("A46").Value = "SOA_formation_24062013_faster"
$ws.Range("B46").Value = "40sect"
$ws.Range("C46").Value = "20day"
$ws.Range("H46").Value = "vector"
$ws.Range("I46").Value = "vap_wallsink"
$ws.Range("J46").Value = "J"

$ws.Range("A47").Value = "temp_20130624T165623"
$ws.Range("B47").Value = 1
$ws.Range("C47").Value = "100nm"
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = "1/(24*3600)"
$ws.Range("F47").Value = 0.3
$ws.Range("G47").Value = 100
$ws.Range("H47").Formula = "=F47*0.0000000000000002*0.00000006*26908000000000000000*G47*0.000000001*26908000000000000000"
$ws.Range("I47").NumberFormat = "d-mmm"
$ws.Range("I47").Value = "1/9s"
$ws.Range("J47").Value = "3nm 1#/cm3s klo 11-13"
$ws.Range("K47").Value = "mass conserv error "

$ws.Range("A48").Value = "temp_20130625T045024"
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = "100nm"
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0.3
$ws.Range("G48").Value = 100
$ws.Range("H48").Formula = "=F48*0.0000000000000002*0.00000006*26908000000000000000*G48*0.000000001*26908000000000000000"
$ws.Range("I48").NumberFormat = "d-mmm"
$ws.Range("I48").Value = "1/9s"
$ws.Range("J48").Value = "3nm 1#/cm3s klo 11-13"
$ws.Range("K48").Value = "mass conserv error "

# ---------------------------------------------------------------------
# Block 2: rows 50-52 (SOA_formation_25062013_test)
# ---------------------------------------------------------------------
$ws.Range("A50").Value = "Script"
$ws.Range("B50").Value = "10sect"
$ws.Range("C50").Value = "20day"
$ws.Range("H50").Font.Bold = $true
$ws.Range("H50").NumberFormat = "#,##0.00"
$ws.Range("H50").Value = "constant"

$ws.Range("A51").Value = "SOA_formation_25062013_test"
$ws.Range("B51").Value = 1
$ws.Range("C51").Value = "100nm"
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0.3
$ws.Range("G51").Value = 0.25
$ws.Range("H51").Formula = "=F51*0.0000000000000002*0.00000006*26908000000000000000*G51*0.000000001*26908000000000000000"
$ws.Range("I51").Value = "1/9000s"
$ws.Range("J51").Value = "3nm 1#/cm3s klo 11-13"
$ws.Range("K51").Value = "mass conserv error "

$ws.Range("A52").Value = "run_20130625T130812"

# ---------------------------------------------------------------------
# Block 3: rows 54-55 (EHDOTUS) - styled with "Explanatory Text"
# ---------------------------------------------------------------------
$ws.Range("A54:J54").Style = "Explanatory Text"
$ws.Range("A54").Value = "Script"
$ws.Range("C54").Value = "1day"
$ws.Range("H54").Value = "exponential"

$ws.Range("A55:J55").Style = "Explanatory Text"
$ws.Range("A55").Value = "EHDOTUS"
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = "100nm"
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0.3
$ws.Range("G55").Value = 0.25
$ws.Range("I55").Value = "1/9000s"
$ws.Range("J55").Value = "3nm 1#/cm3s klo 11-13"

# ---------------------------------------------------------------------
# Block 4: rows 57-59 (SOA_formation_25062013_test2 / used chamber_runfile2)
# ---------------------------------------------------------------------
$ws.Range("A57").Value = "Script"
$ws.Range("B57").Value = "10sect"
$ws.Range("C57").Value = "20day"
$ws.Range("H57").Value = "constant"
$ws.Range("J57").Font.Bold = $true
$ws.Range("J57").Value = "constant"

$ws.Range("A58").Value = "SOA_formation_25062013_test2"
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = "100nm"
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0.3
$ws.Range("G58").Value = 0.25
$ws.Range("H58").Formula = "=F58*0.0000000000000002*0.00000006*26908000000000000000*G58*0.000000001*26908000000000000000"
$ws.Range("I58").Value = "1/9000s"
$ws.Range("J58").Value = "3nm 1#/cm3s "
$ws.Range("K58").Value = "used chamber_runfile2"

$ws.Range("A59").Value = "run_20130625T134834"
$ws.Range("K59").Value = "mass conserv error "

# ---------------------------------------------------------------------
# Block 5: rows 61-63 (SOA_formation_25062013_test2 / used chamber_runfile)
# ---------------------------------------------------------------------
$ws.Range("A61").Value = "Script"
$ws.Range("B61").Value = "10sect"
$ws.Range("C61").Value = "20day"
$ws.Range("H61").Value = "constant"
$ws.Range("J61").Font.Bold = $true
$ws.Range("J61").Value = "constant"

$ws.Range("A62").Value = "SOA_formation_25062013_test2"
$ws.Range("B62").Value = 1
$ws.Range("C62").Value = "100nm"
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0.3
$ws.Range("G62").Value = 0.25
$ws.Range("H62").Formula = "=F62*0.0000000000000002*0.00000006*26908000000000000000*G62*0.000000001*26908000000000000000"
$ws.Range("I62").Value = "1/9000s"
$ws.Range("J62").Value = "3nm 1#/cm3s "
$ws.Range("K62").Value = "used chamber_runfile"

$ws.Range("A63").Value = "run_20130625T141155"
$ws.Range("K63").Value = "mass conserv error "

# ---------------------------------------------------------------------
# Block 6: rows 65-67 (SOA_formation_25062013_test3 / Ntot initial virheellinen)
# ---------------------------------------------------------------------
$ws.Range("A65").Value = "Script"
$ws.Range("B65").Value = "10sect"
$ws.Range("C65").Value = "20day"
$ws.Range("H65").Value = "constant"
$ws.Range("J65").Font.Bold = $true
$ws.Range("J65").Value = "constant"

$ws.Range("A66").Value = "SOA_formation_25062013_test3"
$ws.Range("B66").Value = 1
$ws.Range("C66").Value = "100nm"
$ws.Range("D66").Value = 1000
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 0.3
$ws.Range("G66").Value = 0.25
$ws.Range("H66").Formula = "=F66*0.0000000000000002*0.00000006*26908000000000000000*G66*0.000000001*26908000000000000000"
$ws.Range("I66").Value = "1/9000s"
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = "ei tullut mass conserv error"

$ws.Range("A67").Value = "run_20130625T143638"
$ws.Range("K67").Value = "Ntot initial virheellinen"

# ---------------------------------------------------------------------
# Block 7: rows 69-71 (SOA_formation_25062013_test4 / Ntot initial oikein)
# ---------------------------------------------------------------------
$ws.Range("A69").Value = "Script"
$ws.Range("B69").Value = "20sect"
$ws.Range("C69").Value = "20day"
$ws.Range("H69").Value = "constant"
$ws.Range("J69").Font.Bold = $true
$ws.Range("J69").Value = "constant"

$ws.Range("A70").Value = "SOA_formation_25062013_test4"
$ws.Range("B70").Value = 1
$ws.Range("C70").Value = "100nm"
$ws.Range("D70").Value = 1000
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0.3
$ws.Range("G70").Value = 0.25
$ws.Range("H70").Formula = "=F70*0.0000000000000002*0.00000006*26908000000000000000*G70*0.000000001*26908000000000000000"
$ws.Range("I70").Value = "1/9000s"
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = "ei tullut mass conserv error"

$ws.Range("A71").Value = "run_20130625T145621"
$ws.Range("K71").Value = "Ntot initial oikein"

# ---------------------------------------------------------------------
# Block 8: rows 74-76 (SOA_formation_25062013 / 80 ajoa)
# ---------------------------------------------------------------------
$ws.Range("A74").Value = "Script"

$ws.Range("A75").Value = "SOA_formation_25062013"
$ws.Range("B75").Font.Bold = $true
$ws.Range("B75").Value = "40sect"
$ws.Range("C75").Font.Bold = $true
$ws.Range("C75").Value = "20day"
$ws.Range("D75").Font.Bold = $true
$ws.Range("D75").Value = "80 ajoa"
$ws.Range("H75").Value = "vector"
$ws.Range("I75").Value = "vap_wallsink"
$ws.Range("J75").Value = "J"

$ws.Range("B76").Value = 1
$ws.Range("C76").Value = "10nm"
$ws.Range("D76").Value = 2000
$ws.Range("E76").Value = "1/(24*3600)"
$ws.Range("F76").Value = 0.3
$ws.Range("G76").Value = 100
$ws.Range("H76").Formula = "=F76*0.0000000000000002*0.00000006*26908000000000000000*G76*0.000000001*26908000000000000000"
$ws.Range("I76").NumberFormat = "d-mmm"
$ws.Range("I76").Value = "1/9s"
$ws.Range("J76").Value = 0

# ---------------------------------------------------------------------
# Update the view selection to match the final state
# ---------------------------------------------------------------------
$ws.Range("F82").Select()
